$wb = $excel.ActiveWorkbook

# --- Sheet: contact_data ---
$ws = $wb.Worksheets.Item("contact_data")
$ws.Cells.Item(17, 1).Value = 1
$ws.Cells.Item(17, 2).Value = "2024-12-09 09:10:50"
$ws.Cells.Item(17, 5).Value = "SKZ"

$ws.Cells.Item(18, 1).Value = 1
$ws.Cells.Item(18, 2).Value = "2024-12-09 16:45:18"

$ws.Cells.Item(19, 1).Value = 1
$ws.Cells.Item(19, 2).Value = "2024-12-09 16:45:46"
$ws.Cells.Item(19, 3).Value = "Stefan"
$ws.Cells.Item(19, 4).Value = "Trieß"
$ws.Cells.Item(19, 5).Value = "SKZ"
$ws.Cells.Item(19, 6).Value = "Scientist"
$ws.Cells.Item(19, 7).Value = "s.triess@skz.de"
# "+49" looks numeric to Excel's auto-detection; a leading apostrophe keeps it text,
# matching the source data (inline string "+49", not the number 49).
$ws.Cells.Item(19, 8).Value = "'+49"

# --- Sheet: company_data ---
$ws = $wb.Worksheets.Item("company_data")
$ws.Cells.Item(19, 1).Value = 1
$ws.Cells.Item(19, 2).Value = "2024-12-09 09:11:34"
$ws.Cells.Item(19, 3).Value = "SKZ"
$ws.Cells.Item(19, 4).Value = "Friedrich-Bergius-Ring 22"
# "97076" is a postal code stored as text in the source data.
$ws.Cells.Item(19, 5).Value = "'97076"
$ws.Cells.Item(19, 6).Value = "Würzburg"
$ws.Cells.Item(19, 7).Value = "Bayern"
$ws.Cells.Item(19, 8).Value = "Deutschland"
$ws.Cells.Item(19, 9).Value = $false
$ws.Cells.Item(19, 10).Value = $false
$ws.Cells.Item(19, 12).Value = 49.80282025
$ws.Cells.Item(19, 13).Value = 10.00010726291456

$ws.Cells.Item(20, 1).Value = 1
$ws.Cells.Item(20, 2).Value = "2024-12-09 16:46:07"
$ws.Cells.Item(20, 3).Value = "SKZ"
$ws.Cells.Item(20, 4).Value = "Friedrich-Bergius-Ring 22"
$ws.Cells.Item(20, 5).Value = "'97076"
$ws.Cells.Item(20, 6).Value = "Würzburg"
$ws.Cells.Item(20, 7).Value = "Bayern"
$ws.Cells.Item(20, 8).Value = "Deutschland"
$ws.Cells.Item(20, 9).Value = $false
$ws.Cells.Item(20, 10).Value = $false
$ws.Cells.Item(20, 12).Value = 49.80282025
$ws.Cells.Item(20, 13).Value = 10.00010726291456

# --- Sheet: product_fractions ---
$ws = $wb.Worksheets.Item("product_fractions")
$ws.Cells.Item(17, 1).Value = 1
$ws.Cells.Item(17, 2).Value = "2024-12-09 09:11:53"
$ws.Cells.Item(17, 3).Value = "['PE-LD', 'PP', 'ABS', 'Magnesium']"
$ws.Cells.Item(17, 4).Value = "['', '', '', '']"
$ws.Cells.Item(17, 5).Value = "[80.0, 10.0, 5.0, 5.0]"

$ws.Cells.Item(18, 1).Value = 1
$ws.Cells.Item(18, 2).Value = "2024-12-09 09:28:05"
$ws.Cells.Item(18, 3).Value = "['PE-LD', 'PE-MD', 'ABS', 'Magnesium']"
$ws.Cells.Item(18, 4).Value = "['', '', '', '']"
$ws.Cells.Item(18, 5).Value = "[80.0, 10.0, 5.0, 5.0]"

$ws.Cells.Item(19, 1).Value = 1
$ws.Cells.Item(19, 2).Value = "2024-12-09 09:28:23"
$ws.Cells.Item(19, 3).Value = "['PE-LD', 'PP', 'ABS', 'Magnesium']"
$ws.Cells.Item(19, 4).Value = "['', '', '', '']"
$ws.Cells.Item(19, 5).Value = "[80.0, 10.0, 5.0, 5.0]"

$ws.Cells.Item(20, 1).Value = 1
$ws.Cells.Item(20, 2).Value = "2024-12-09 10:11:50"
$ws.Cells.Item(20, 3).Value = "['PE-LD', 'PS', 'ABS', 'Magnesium']"
$ws.Cells.Item(20, 4).Value = "['', '', '', '']"
$ws.Cells.Item(20, 5).Value = "[80.0, 10.0, 5.0, 5.0]"

$ws.Cells.Item(21, 1).Value = 1
$ws.Cells.Item(21, 2).Value = "2024-12-09 16:46:32"
$ws.Cells.Item(21, 3).Value = "['PS', 'PET', 'PA', 'Eisen']"
$ws.Cells.Item(21, 4).Value = "['', '', '', '']"
$ws.Cells.Item(21, 5).Value = "[60.0, 20.0, 10.0, 10.0]"

# --- Sheet: product_origin ---
$ws = $wb.Worksheets.Item("product_origin")
$ws.Cells.Item(9, 1).Value = 1
$ws.Cells.Item(9, 2).Value = "2024-12-09 09:11:59"
$ws.Cells.Item(9, 3).Value = "Post-Industrial (PI)"

$ws.Cells.Item(10, 1).Value = 1
$ws.Cells.Item(10, 2).Value = "2024-12-09 16:46:37"
$ws.Cells.Item(10, 3).Value = "Post-Industrial (PI)"

# --- Sheet: product_quality ---
$ws = $wb.Worksheets.Item("product_quality")
$ws.Cells.Item(14, 1).Value = 1
$ws.Cells.Item(14, 2).Value = "2024-12-09 09:12:50"
$ws.Cells.Item(14, 3).Value = "Ja"
$ws.Cells.Item(14, 4).Value = "blau"
$ws.Cells.Item(14, 5).Value = 100
$ws.Cells.Item(14, 6).Value = "mittel"
$ws.Cells.Item(14, 8).Value = "[[], [], [], []]"
$ws.Cells.Item(14, 9).Value = "[[], [], [], []]"

$ws.Cells.Item(15, 1).Value = 1
$ws.Cells.Item(15, 2).Value = "2024-12-09 16:46:50"
$ws.Cells.Item(15, 3).Value = "Ja"
$ws.Cells.Item(15, 4).Value = "gelb"
$ws.Cells.Item(15, 5).Value = 100
$ws.Cells.Item(15, 6).Value = "mittel"
$ws.Cells.Item(15, 8).Value = "[[], [], [], []]"
$ws.Cells.Item(15, 9).Value = "[[], [], [], []]"

$ws.Cells.Item(16, 1).Value = 1
$ws.Cells.Item(16, 2).Value = "2024-12-09 17:01:04"
$ws.Cells.Item(16, 3).Value = "Ja"
$ws.Cells.Item(16, 4).Value = "gelb"
$ws.Cells.Item(16, 5).Value = 100
$ws.Cells.Item(16, 6).Value = "sehr hoch"
$ws.Cells.Item(16, 8).Value = "[[], [], [], []]"
$ws.Cells.Item(16, 9).Value = "[[], [], [], []]"

$ws.Cells.Item(17, 1).Value = 1
$ws.Cells.Item(17, 2).Value = "2024-12-09 17:03:09"
$ws.Cells.Item(17, 3).Value = "Ja"
$ws.Cells.Item(17, 4).Value = "gelb"
$ws.Cells.Item(17, 5).Value = 100
$ws.Cells.Item(17, 6).Value = "keine"
$ws.Cells.Item(17, 8).Value = "[[], [], [], []]"
$ws.Cells.Item(17, 9).Value = "[[], [], [], []]"

$ws.Cells.Item(18, 1).Value = 1
$ws.Cells.Item(18, 2).Value = "2024-12-09 17:03:34"
$ws.Cells.Item(18, 3).Value = "Ja"
$ws.Cells.Item(18, 4).Value = "gelb"
$ws.Cells.Item(18, 5).Value = 100
$ws.Cells.Item(18, 6).Value = "gering"
$ws.Cells.Item(18, 8).Value = "[[], [], [], []]"
$ws.Cells.Item(18, 9).Value = "[[], [], [], []]"

$ws.Cells.Item(19, 1).Value = 1
$ws.Cells.Item(19, 2).Value = "2024-12-09 17:04:04"
$ws.Cells.Item(19, 3).Value = "Ja"
$ws.Cells.Item(19, 4).Value = "gelb"
$ws.Cells.Item(19, 5).Value = 100
$ws.Cells.Item(19, 6).Value = "mittel"
# source row has an explicit-but-empty "Verschmutzungsart" cell (form field left blank
# but submitted); a lone apostrophe reproduces a present empty text cell instead of
# leaving the cell absent.
$ws.Cells.Item(19, 7).Value = "'"
$ws.Cells.Item(19, 8).Value = "[[], [], [], []]"
$ws.Cells.Item(19, 9).Value = "[[], [], [], []]"

# --- Sheet: product_amount ---
$ws = $wb.Worksheets.Item("product_amount")
# Row 17's "Weitere Beschreibung" cell was an empty placeholder left over from an
# earlier (blank) form submission; the refreshed export no longer carries it.
$ws.Cells.Item(17, 6).ClearContents()

$ws.Cells.Item(18, 1).Value = 1
$ws.Cells.Item(18, 2).Value = "2024-12-09 09:13:02"
$ws.Cells.Item(18, 3).Value = 5
$ws.Cells.Item(18, 4).Value = 5
$ws.Cells.Item(18, 5).Value = "Woche"

$ws.Cells.Item(19, 1).Value = 1
$ws.Cells.Item(19, 2).Value = "2024-12-09 16:47:02"
$ws.Cells.Item(19, 3).Value = 5
$ws.Cells.Item(19, 4).Value = 10
$ws.Cells.Item(19, 5).Value = "Monat"
